$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.234.14"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.645.81"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'217.25"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "'0.0639"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'20.00"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "'4.31"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.873.35"
$ws.Range("D14").Value = "1.626.81"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'0.550"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "'63.56"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "26.214.71"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'4.44"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'195.34"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "'10.06"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'6.35"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'143.38"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("D28").Value = "'6.94"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'15.62"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").Value = "'0.0505"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'1.61"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").Value = "'0.914"
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.555"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.133.96"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").Value = "'0.0157"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'5.66"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "'100.18"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'0.797"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "1.782.57"
$ws.Range("D46").Value = "'56.33"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.48"
$ws.Range("E47").Value = "  +4.38%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0516"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.418"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.0976"
$ws.Range("E51").Value = "  +2.27%  "